$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

try {
  $ws.Range("E15:E28").Validation.Delete()
  Write-Output "deleted ok"
} catch {
  Write-Output "delete err: $_"
}
